# Edit the delivery-application "form.xlsx" weighing sheet:
#  - move the active selection
#  - add a TODAY() formula for the date field
#  - wire up the per-class storage / net-weight / payable calculations
#  - wire up the totals / averages section at the bottom

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Selection moved from A11 to E8 --------------------------------------
$ws.Range("E8").Select()

# --- Date field: B6 now computes today's date -----------------------------
$ws.Range("B6").Formula = "=TODAY()"

# --- Per-class rows (25-30): Bruttó kg already filled in, fill in the
#     "levonás" (deduction) input with 0 and wire up Nettó kg / Kifizetendő
#     bruttó formulas -----------------------------------------------------
$ws.Range("C25").Value = 0
$ws.Range("C26").Value = 0
$ws.Range("C27").Value = 0
$ws.Range("C28").Value = 0
$ws.Range("C29").Value = 0
$ws.Range("C30").Value = 0

$ws.Range("E25").Formula = "=ROUNDDOWN(C25*0.95,0)"
$ws.Range("E26").Formula = "=ROUNDDOWN(C26*0.95,0)"
$ws.Range("E27").Formula = "=ROUNDDOWN(C27*0.95,0)"
$ws.Range("E28").Value = 0
$ws.Range("E29").Value = 0
$ws.Range("E30").Formula = "=ROUNDDOWN(C30*0.95,0)"

$ws.Range("F25:F30").Formula = "=SUM(B25*E25)"

# --- Totals row (31) and "előleg" row (33) --------------------------------
$ws.Range("E31").Formula = "=SUM(E25:E30)"
$ws.Range("F31").Formula = "=SUM(F25:F30)"
$ws.Range("F33").Formula = "=SUM(F31-F32)"

# --- Summary (40-41): Átlagár / Kifizetendő section -----------------------
$ws.Range("C40").Formula = "=ROUND(F40/1.12,0)"
$ws.Range("D40").Formula = "=SUM(F40-C40)"
$ws.Range("E40").Formula = "=ROUND(F31/E31,2)"
$ws.Range("F40").Formula = "=ROUND(E31*E40,0)"
$ws.Range("E41").Formula = "=ROUND(E40/1.12,2)"

$wb.Save()
